$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 210307
$ws.Range("I19").Value = 0
$ws.Range("J19").Value = 210307
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 210307
$ws.Range("M19").ClearContents()
$ws.Range("N19").Value = -210657

$ws.Range("H33").Value = 15627053
$ws.Range("I33").Value = 25001536
$ws.Range("J33").Value = 2916.6667
$ws.Range("K33").Value = 25001536
$ws.Range("L33").Value = 2916.6667
$ws.Range("M33").Value = -25001307
$ws.Range("N33").Value = -3374.6667

$ws.Range("H40").Value = 1107.5
$ws.Range("I40").Value = 1015
$ws.Range("J40").Value = 1200
$ws.Range("K40").Value = 1015
$ws.Range("L40").Value = 1200
$ws.Range("M40").Value = -840
$ws.Range("N40").Value = -1550

$ws.Range("H86").Value = 13245.667
$ws.Range("I86").Value = 12993.875
$ws.Range("J86").Value = 13749.25
$ws.Range("K86").Value = 12993.875
$ws.Range("L86").Value = 13749.25
$ws.Range("M86").Value = -11870.875
$ws.Range("N86").Value = -15995.25

$ws.Range("H89").Value = 13245.667
$ws.Range("I89").Value = 12993.875
$ws.Range("J89").Value = 13749.25
$ws.Range("K89").Value = 64969.375
$ws.Range("L89").Value = 68746.25
$ws.Range("M89").Value = -59353.375
$ws.Range("N89").Value = -79978.25

$ws.Range("H125").Value = 3797998.8
$ws.Range("J125").Value = 10588
$ws.Range("L125").Value = 95292
$ws.Range("N125").Value = -100212

$ws.Range("H137").Value = 14686.4375
$ws.Range("I137").Value = 1487.4
$ws.Range("K137").Value = 4462.200000000001
$ws.Range("M137").Value = -1912.200000000001

$ws.Range("H138").Value = 3628.373
$ws.Range("I138").Value = 2763.4167
$ws.Range("J138").Value = 3817.0908
$ws.Range("K138").Value = 8290.250100000001
$ws.Range("L138").Value = 11451.2724
$ws.Range("M138").Value = -3150.250100000001
$ws.Range("N138").Value = -21731.2724

$ws.Range("H141").Value = 3280.4333
$ws.Range("I141").Value = 3119.7083
$ws.Range("J141").Value = 3923.3333
$ws.Range("K141").Value = 9359.124899999999
$ws.Range("L141").Value = 11769.9999
$ws.Range("M141").Value = -4179.124899999999
$ws.Range("N141").Value = -22129.9999


$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("N39").ClearContents()

$ws.Range("H61").Value = 735266.2
$ws.Range("I61").Value = 2412.946
$ws.Range("K61").Value = 2412.946
$ws.Range("M61").Value = -2200.946

$ws.Range("H122").Value = 3459825.5
$ws.Range("I122").Value = 6912651.5
$ws.Range("K122").Value = 20737954.5
$ws.Range("M122").Value = -20735504.5

$ws.Range("H136").Value = 735266.2
$ws.Range("I136").Value = 2412.946
$ws.Range("K136").Value = 7238.838
$ws.Range("M136").Value = -4688.838


$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 1312.5
$ws.Range("I11").Value = 307.5
$ws.Range("J11").Value = 2317.5
$ws.Range("K11").Value = 307.5
$ws.Range("L11").Value = 2317.5
$ws.Range("M11").Value = -167.5
$ws.Range("N11").Value = -2597.5

$ws.Range("H20").Value = 7770530.5
$ws.Range("I20").Value = 16672408
$ws.Range("K20").Value = 16672408
$ws.Range("M20").Value = -16672161

$ws.Range("H46").Value = 4500
$ws.Range("J46").Value = 4500
$ws.Range("L46").Value = 4500
$ws.Range("N46").Value = -5096

$ws.Range("H55").Value = 149329.67
$ws.Range("I55").Value = 98994
$ws.Range("K55").Value = 98994
$ws.Range("M55").Value = -98721

$ws.Range("H82").Value = 16889.5
$ws.Range("I82").Value = 7519.3335
$ws.Range("K82").Value = 7519.3335
$ws.Range("M82").Value = -7136.3335

$ws.Range("H85").Value = 16889.5
$ws.Range("I85").Value = 7519.3335
$ws.Range("K85").Value = 7519.3335
$ws.Range("M85").Value = -6193.3335

$ws.Range("H86").Value = 111114760
$ws.Range("I86").Value = 4218.5713
$ws.Range("K86").Value = 4218.5713
$ws.Range("M86").Value = -3095.5713

$ws.Range("H89").Value = 111114760
$ws.Range("I89").Value = 4218.5713
$ws.Range("K89").Value = 21092.8565
$ws.Range("M89").Value = -15476.8565

$ws.Range("H94").Value = 638.7646999999999
$ws.Range("I94").Value = 726.12
$ws.Range("K94").Value = 726.12
$ws.Range("M94").Value = -275.12

$ws.Range("H135").Value = 84614.58
$ws.Range("J135").Value = 84614.58
$ws.Range("L135").Value = 84614.58
$ws.Range("N135").Value = -94754.58


$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6250.6587
$ws.Range("I31").Value = 1201.5758
$ws.Range("J31").Value = 27078.125
$ws.Range("K31").Value = 1201.5758
$ws.Range("L31").Value = 27078.125
$ws.Range("M31").Value = -906.5758000000001
$ws.Range("N31").Value = -27668.125

$ws.Range("H34").Value = 6250.6587
$ws.Range("I34").Value = 1201.5758
$ws.Range("J34").Value = 27078.125
$ws.Range("K34").Value = 1201.5758
$ws.Range("L34").Value = 27078.125
$ws.Range("M34").Value = -999.5758000000001
$ws.Range("N34").Value = -27482.125

$ws.Range("H62").Value = 3004.5
$ws.Range("I62").Value = 3
$ws.Range("J62").Value = 6006
$ws.Range("K62").Value = 3
$ws.Range("L62").Value = 6006
$ws.Range("M62").Value = 621
$ws.Range("N62").Value = -7254

$ws.Range("H65").Value = 3004.5
$ws.Range("I65").Value = 3
$ws.Range("J65").Value = 6006
$ws.Range("K65").Value = 15
$ws.Range("L65").Value = 30030
$ws.Range("M65").Value = 3105
$ws.Range("N65").Value = -36270


$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 3906.4
$ws.Range("J39").Value = 4139.3125
$ws.Range("L39").Value = 12417.9375
$ws.Range("N39").Value = -13005.9375

$ws.Range("H56").Value = 7162.5
$ws.Range("I56").Value = 7162.5
$ws.Range("K56").Value = 7162.5
$ws.Range("M56").Value = -6632.5

$ws.Range("H68").Value = 2317.2144
$ws.Range("J68").Value = 1480.5834
$ws.Range("L68").Value = 4441.7502
$ws.Range("N68").Value = -6063.7502

$ws.Range("H71").Value = 2317.2144
$ws.Range("J71").Value = 1480.5834
$ws.Range("L71").Value = 13325.2506
$ws.Range("N71").Value = -21437.2506

$ws.Range("H80").Value = 35000
$ws.Range("I80").Value = 20000
$ws.Range("K80").Value = 60000
$ws.Range("M80").Value = -59064

$ws.Range("H83").Value = 35000
$ws.Range("I83").Value = 20000
$ws.Range("K83").Value = 180000
$ws.Range("M83").Value = -175320

$ws.Range("H131").Value = 1376.34
$ws.Range("J131").Value = 1499.9878
$ws.Range("L131").Value = 4499.963400000001
$ws.Range("N131").Value = -14579.9634


$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 112985
$ws.Range("J42").Value = 112985
$ws.Range("L42").Value = 112985
$ws.Range("N42").Value = -113955

$ws.Range("H70").Value = 6096
$ws.Range("I70").Value = 4716
$ws.Range("J70").Value = 7016
$ws.Range("K70").Value = 4716
$ws.Range("L70").Value = 7016
$ws.Range("M70").Value = -4446
$ws.Range("N70").Value = -7556

$ws.Range("H73").Value = 6096
$ws.Range("I73").Value = 4716
$ws.Range("J73").Value = 7016
$ws.Range("K73").Value = 4716
$ws.Range("L73").Value = 7016
$ws.Range("M73").Value = -3780
$ws.Range("N73").Value = -8888

$ws.Range("H80").Value = 3899.923
$ws.Range("I80").Value = 3936.75
$ws.Range("J80").Value = 3841
$ws.Range("K80").Value = 3936.75
$ws.Range("L80").Value = 3841
$ws.Range("M80").Value = -2938.75
$ws.Range("N80").Value = -5837

$ws.Range("H83").Value = 3899.923
$ws.Range("I83").Value = 3936.75
$ws.Range("J83").Value = 3841
$ws.Range("K83").Value = 19683.75
$ws.Range("L83").Value = 19205
$ws.Range("M83").Value = -14691.75
$ws.Range("N83").Value = -29189

$ws.Range("H92").Value = 9951.833000000001
$ws.Range("J92").Value = 9951.833000000001
$ws.Range("L92").Value = 9951.833000000001
$ws.Range("N92").Value = -13695.833

$ws.Range("H95").Value = 29000
$ws.Range("J95").Value = 29000
$ws.Range("L95").Value = 29000
$ws.Range("N95").Value = -34492

$ws.Range("H102").Value = 4099193
$ws.Range("I102").Value = 5009470.5
$ws.Range("K102").Value = 5009470.5
$ws.Range("M102").Value = -5007848.5

$ws.Range("H115").Value = 112985
$ws.Range("J115").Value = 112985
$ws.Range("L115").Value = 112985
$ws.Range("N115").Value = -115335

$ws.Range("H126").Value = 7671432.5
$ws.Range("I126").Value = 3595498.5
$ws.Range("J126").Value = 17861268
$ws.Range("K126").Value = 10786495.5
$ws.Range("L126").Value = 53583804
$ws.Range("M126").Value = -10784025.5
$ws.Range("N126").Value = -53588744

$ws.Range("H141").Value = 200000
$ws.Range("J141").Value = 200000
$ws.Range("L141").Value = 200000
$ws.Range("N141").Value = -210360


$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 45456148
$ws.Range("I16").Value = 47620690
$ws.Range("J16").Value = 808
$ws.Range("K16").Value = 47620690
$ws.Range("L16").Value = 808
$ws.Range("M16").Value = -47620520
$ws.Range("N16").Value = -1148

$ws.Range("H22").Value = 58824864
$ws.Range("J22").Value = 125001150
$ws.Range("L22").Value = 125001150
$ws.Range("N22").Value = -125001740

$ws.Range("H27").Value = 58824864
$ws.Range("J27").Value = 125001150
$ws.Range("L27").Value = 125001150
$ws.Range("N27").Value = -125001364

$ws.Range("H40").Value = 58823530
$ws.Range("I40").Value = 0
$ws.Range("J40").Value = 58823530
$ws.Range("K40").Value = 0
$ws.Range("L40").Value = 58823530
$ws.Range("M40").ClearContents()
$ws.Range("N40").Value = -58823802

$ws.Range("H82").Value = 1934.4445
$ws.Range("I82").Value = 2233.7144
$ws.Range("K82").Value = 2233.7144
$ws.Range("M82").Value = -1872.7144

$ws.Range("H85").Value = 1934.4445
$ws.Range("I85").Value = 2233.7144
$ws.Range("K85").Value = 2233.7144
$ws.Range("M85").Value = -985.7143999999998

$ws.Range("H93").Value = 100005950
$ws.Range("I93").Value = 125007130
$ws.Range("J93").Value = 1249.5
$ws.Range("K93").Value = 125007130
$ws.Range("L93").Value = 1249.5
$ws.Range("M93").Value = -125005882
$ws.Range("N93").Value = -3745.5

$ws.Range("H100").Value = 2513.6155
$ws.Range("J100").Value = 2844.5
$ws.Range("L100").Value = 2844.5
$ws.Range("N100").Value = -3926.5

$ws.Range("H104").Value = 28333
$ws.Range("J104").Value = 28333
$ws.Range("L104").Value = 28333
$ws.Range("N104").Value = -35321

$ws.Range("H122").Value = 41986690
$ws.Range("I122").Value = 49601256
$ws.Range("K122").Value = 148803768
$ws.Range("M122").Value = -148801318

$ws.Range("H132").Value = 2796719.2
$ws.Range("I132").Value = 4718
$ws.Range("J132").Value = 4990434.5
$ws.Range("K132").Value = 14154
$ws.Range("L132").Value = 14971303.5
$ws.Range("M132").Value = -11624
$ws.Range("N132").Value = -14976363.5


$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 793.3684
$ws.Range("I107").Value = 986.53845
$ws.Range("K107").Value = 2959.61535
$ws.Range("M107").Value = -1039.61535

$ws.Range("H119").Value = 226666.67
$ws.Range("J119").Value = 226666.67
$ws.Range("L119").Value = 226666.67
$ws.Range("N119").Value = -236342.67

$ws.Range("H122").Value = 859094.8
$ws.Range("I122").Value = 1588854.6
$ws.Range("K122").Value = 4766563.800000001
$ws.Range("M122").Value = -4764113.800000001

$ws.Range("H123").Value = 85437.5
$ws.Range("J123").Value = 85437.5
$ws.Range("L123").Value = 85437.5
$ws.Range("N123").Value = -95237.5

$ws.Range("H125").Value = 50000
$ws.Range("J125").Value = 50000
$ws.Range("L125").Value = 50000
$ws.Range("N125").Value = -59840

$ws.Range("H126").Value = 5019753
$ws.Range("I126").Value = 23113.354
$ws.Range("K126").Value = 69340.06200000001
$ws.Range("M126").Value = -66870.06200000001

